$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 0.0467
$ws.Range("E2").Value = 0.0984
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 0.1088491898628999
$ws.Range("H2").Value = 0.1088491898628999
$ws.Range("I2").Value = 0.09056917324470296
$ws.Range("J2").Value = 0.06896882617439548
$ws.Range("K2").Value = 17.3
$ws.Range("L2").Value = 0.0718737017033652
$ws.Range("M2").Value = 10.63
$ws.Range("N2").Value = 0.06796675191815855
$ws.Range("O2").Value = 0.614450867052023
$ws.Range("P2").Value = 8.02
$ws.Range("Q2").Value = 0.05127877237851662
$ws.Range("R2").Value = 0.4635838150289017
$ws.Range("S2").Value = 2.609999999999999
$ws.Range("T2").Value = 0.2455315145813734
$ws.Range("U2").Value = 9.49
$ws.Range("V2").Value = 0.06067774936061381
$ws.Range("W2").Value = 0.1146454605699139
$ws.Range("X2").Value = 0.08899879257416607
$ws.Range("Y2").Value = 0.02564666799574779
$ws.Range("Z2").Value = 1.693043539424633
$ws.Range("AA2").Value = 0.1167672255762608
$ws.Range("AB2").Value = 0.08855629968098354
$ws.Range("AC2").Value = 0.02821092589527724
$ws.Range("AD2").Value = 1.39
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1.39
$ws.Range("AG2").Value = -8.1
$ws.Range("AH2").Value = 0.008809176753913428
$ws.Range("AI2").Value = 0.008759216081668662
$ws.Range("AJ2").Value = -0.05461901550910316
$ws.Range("AK2").Value = -0.05428954423592493
$ws.Range("AL2").Value = 0.024
$ws.Range("AM2").Value = 0.024
$ws.Range("AN2").Value = 0.0594017094017094
$ws.Range("AO2").Value = 908.3333333333334
$ws.Range("AP2").Value = -0.3461538461538461
$ws.Range("AQ2").Value = 908.3333333333334

# Row 3 updates
$ws.Range("D3").Value = 0.0467
$ws.Range("E3").Value = 0.0984
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = 0.1088491898628999
$ws.Range("H3").Value = 0.1088491898628999
$ws.Range("I3").Value = 0.09056917324470296
$ws.Range("J3").Value = 0.06896882617439548
$ws.Range("K3").Value = 17.3
$ws.Range("L3").Value = 0.0718737017033652
$ws.Range("M3").Value = 10.63
$ws.Range("N3").Value = 0.06796675191815855
$ws.Range("O3").Value = 0.614450867052023
$ws.Range("P3").Value = 8.02
$ws.Range("Q3").Value = 0.05127877237851662
$ws.Range("R3").Value = 0.4635838150289017
$ws.Range("S3").Value = 2.609999999999999
$ws.Range("T3").Value = 0.2455315145813734
$ws.Range("U3").Value = 9.49
$ws.Range("V3").Value = 0.06067774936061381
$ws.Range("W3").Value = 0.1146454605699139
$ws.Range("X3").Value = 0.08899879257416607
$ws.Range("Y3").Value = 0.02564666799574779
$ws.Range("Z3").Value = 1.693043539424633
$ws.Range("AA3").Value = 0.1167672255762608
$ws.Range("AB3").Value = 0.08855629968098354
$ws.Range("AC3").Value = 0.02821092589527724
$ws.Range("AD3").Value = 1.39
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1.39
$ws.Range("AG3").Value = -8.1
$ws.Range("AH3").Value = 0.008809176753913428
$ws.Range("AI3").Value = 0.008759216081668662
$ws.Range("AJ3").Value = -0.05461901550910316
$ws.Range("AK3").Value = -0.05428954423592493
$ws.Range("AL3").Value = 0.024
$ws.Range("AM3").Value = 0.024
$ws.Range("AN3").Value = 0.0594017094017094
$ws.Range("AO3").Value = 908.3333333333334
$ws.Range("AP3").Value = -0.3461538461538461
$ws.Range("AQ3").Value = 908.3333333333334
